$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price (column D) / Volume(1h) (column E) refresh.
# Values that look like plain numbers (e.g. "7.31", "0.999") must be forced to
# Text so Excel does not silently convert them to numeric cells; the cells
# original style is saved/restored around the write so no formatting changes.

$ws.Range("D2").Value = '49.534.68'
$ws.Range("E2").Value = '  -1.16%  '
$ws.Range("D3").Value = '2.638.02'
$ws.Range("E3").Value = '  +0.13%  '
$style = $ws.Range("D4").Style
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("D4").Style = $style
$ws.Range("E4").Value = '  -0.02%  '
$style = $ws.Range("D5").Style
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '112.51'
$ws.Range("D5").Style = $style
$ws.Range("E5").Value = '  +1.45%  '
$style = $ws.Range("D6").Style
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '323.98'
$ws.Range("D6").Style = $style
$ws.Range("E6").Value = '  -1.26%  '
$ws.Range("E7").Value = '  -1.13%  '
$style = $ws.Range("D8").Style
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.999'
$ws.Range("D8").Style = $style
$ws.Range("E8").Value = '  +0.02%  '
$style = $ws.Range("D9").Style
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.544'
$ws.Range("D9").Style = $style
$ws.Range("E9").Value = '  -3.12%  '
$style = $ws.Range("D10").Style
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '39.92'
$ws.Range("D10").Style = $style
$ws.Range("E10").Value = '  -1.81%  '
$style = $ws.Range("D11").Style
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '19.83'
$ws.Range("D11").Style = $style
$ws.Range("E11").Value = '  -4.19%  '
$style = $ws.Range("D12").Style
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0813'
$ws.Range("D12").Style = $style
$ws.Range("E12").Value = '  -0.97%  '
$style = $ws.Range("D13").Style
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.127'
$ws.Range("D13").Style = $style
$ws.Range("E13").Value = '  +1.37%  '
$style = $ws.Range("D14").Style
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '7.31'
$ws.Range("D14").Style = $style
$ws.Range("E14").Value = '  +0.21%  '
$ws.Range("D15").Value = '3.041.18'
$ws.Range("E15").Value = '  -0.38%  '
$ws.Range("D16").Value = '2.648.19'
$ws.Range("E16").Value = '  +0.69%  '
$style = $ws.Range("D17").Style
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.861'
$ws.Range("D17").Style = $style
$ws.Range("E17").Value = '  -1.94%  '
$ws.Range("D18").Value = '49.392.16'
$ws.Range("E18").Value = '  -1.25%  '
$ws.Range("E19").Value = '  -2.04%  '
$style = $ws.Range("D20").Style
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '12.92'
$ws.Range("D20").Style = $style
$ws.Range("E20").Value = '  -3.16%  '
$style = $ws.Range("D21").Style
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.71'
$ws.Range("D21").Style = $style
$ws.Range("E21").Value = '  -1.88%  '
$ws.Range("D22").Value = '0.0₃0948'
$ws.Range("E22").Value = '  -1.22%  '
$style = $ws.Range("D23").Style
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '270.49'
$ws.Range("D23").Style = $style
$ws.Range("E23").Value = '  -3.35%  '
$style = $ws.Range("D24").Style
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '68.90'
$ws.Range("D24").Style = $style
$ws.Range("E24").Value = '  -5.55%  '
$ws.Range("E25").Value = '  -2.69%  '
$style = $ws.Range("D26").Style
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '26.31'
$ws.Range("D26").Style = $style
$ws.Range("E26").Value = '  -1.22%  '
$ws.Range("E28").Value = '  +3.79%  '
$ws.Range("E29").Value = '  -0.49%  '
$style = $ws.Range("D30").Style
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '35.15'
$ws.Range("D30").Style = $style
$ws.Range("E30").Value = '  -4.76%  '
$ws.Range("E31").Value = '  -4.07%  '
$style = $ws.Range("D32").Style
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '49.59'
$ws.Range("D32").Style = $style
$ws.Range("E32").Value = '  -0.54%  '
$style = $ws.Range("D33").Style
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.50'
$ws.Range("D33").Style = $style
$ws.Range("E33").Value = '  +0.94%  '
$ws.Range("E34").Value = '  +2.53%  '
$ws.Range("E35").Value = '  -0.23%  '
$style = $ws.Range("D36").Style
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '18.96'
$ws.Range("D36").Style = $style
$ws.Range("E36").Value = '  -4.21%  '
$style = $ws.Range("D37").Style
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '4.94'
$ws.Range("D37").Style = $style
$ws.Range("E37").Value = '  +3.80%  '
$style = $ws.Range("D38").Style
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.05'
$ws.Range("D38").Style = $style
$ws.Range("E38").Value = '  -0.58%  '
$style = $ws.Range("D39").Style
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.14'
$ws.Range("D39").Style = $style
$ws.Range("E39").Value = '  +1.08%  '
$style = $ws.Range("D40").Style
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '126.49'
$ws.Range("D40").Style = $style
$ws.Range("E40").Value = '  +2.62%  '
$ws.Range("E41").Value = '  -1.59%  '
$style = $ws.Range("D42").Style
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '22.37'
$ws.Range("D42").Style = $style
$ws.Range("E42").Value = '  -0.90%  '
$style = $ws.Range("D43").Style
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.0324'
$ws.Range("D43").Style = $style
$ws.Range("E43").Value = '  +3.00%  '
$ws.Range("E44").Value = '  -3.69%  '
$ws.Range("D45").Value = '2.061.06'
$ws.Range("E45").Value = '  +0.13%  '
$style = $ws.Range("D46").Style
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.14'
$ws.Range("D46").Style = $style
$ws.Range("E46").Value = '  +6.39%  '
$style = $ws.Range("D47").Style
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.22'
$ws.Range("D47").Style = $style
$ws.Range("E47").Value = '  -4.19%  '
$ws.Range("E48").Value = '  -8.68%  '
$ws.Range("E49").Value = '  -1.19%  '
$style = $ws.Range("D50").Style
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '59.19'
$ws.Range("D50").Style = $style
$ws.Range("E50").Value = '  +0.67%  '
$style = $ws.Range("D51").Style
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '5.22'
$ws.Range("D51").Style = $style
$ws.Range("E51").Value = '  -2.97%  '
